$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet "manual" -> "Beneficiaries"
$ws.Name = "Beneficiaries"

# Update the Gender value in H2 from "Male" -> "M"
$ws.Range("H2").Value = "M"

# Update the active cell selection to I6
$ws.Range("I6").Select()
